# Auto-generated edits applying the scheduled-runner price/profit updates
# to the Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 1000
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1460
$ws.Range("H64").Value = 7324.5
$ws.Range("I64").Value = 4311.5
$ws.Range("J64").Value = 9333.166999999999
$ws.Range("K64").Value = 4311.5
$ws.Range("L64").Value = 9333.166999999999
$ws.Range("M64").Value = -4063.5
$ws.Range("N64").Value = -9829.166999999999
$ws.Range("H67").Value = 7324.5
$ws.Range("I67").Value = 4311.5
$ws.Range("J67").Value = 9333.166999999999
$ws.Range("K67").Value = 4311.5
$ws.Range("L67").Value = 9333.166999999999
$ws.Range("M67").Value = -3453.5
$ws.Range("N67").Value = -11049.167
$ws.Range("H74").Value = 20000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 20000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H80").Value = 473.10345
$ws.Range("I80").Value = 323.0909
$ws.Range("J80").Value = 564.7778
$ws.Range("K80").Value = 969.2727
$ws.Range("L80").Value = 1694.3334
$ws.Range("M80").Value = 28.72730000000001
$ws.Range("N80").Value = -3690.3334
$ws.Range("H83").Value = 473.10345
$ws.Range("I83").Value = 323.0909
$ws.Range("J83").Value = 564.7778
$ws.Range("K83").Value = 2907.8181
$ws.Range("L83").Value = 5083.000199999999
$ws.Range("M83").Value = 2084.1819
$ws.Range("N83").Value = -15067.0002
$ws.Range("H87").Value = 91999.75
$ws.Range("J87").Value = 91999.75
$ws.Range("L87").Value = 91999.75
$ws.Range("N87").Value = -94495.75
$ws.Range("H90").Value = 91999.75
$ws.Range("J90").Value = 91999.75
$ws.Range("L90").Value = 275999.25
$ws.Range("N90").Value = -288479.25
$ws.Range("H105").Value = 34335.5
$ws.Range("J105").Value = 34335.5
$ws.Range("L105").Value = 34335.5
$ws.Range("N105").Value = -41323.5
$ws.Range("H132").Value = 13924.211
$ws.Range("I132").Value = 13272.5
$ws.Range("K132").Value = 39817.5
$ws.Range("M132").Value = -37287.5
$ws.Range("H138").Value = 2851.7646
$ws.Range("J138").Value = 3377.4546
$ws.Range("L138").Value = 10132.3638
$ws.Range("N138").Value = -20412.3638
$ws.Range("H141").Value = 2994.875
$ws.Range("I141").Value = 1342.3334
$ws.Range("K141").Value = 4027.0002
$ws.Range("M141").Value = 1152.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4982.4
$ws.Range("I74").Value = 4982.4
$ws.Range("K74").Value = 4982.4
$ws.Range("M74").Value = -4108.4
$ws.Range("H77").Value = 4982.4
$ws.Range("I77").Value = 4982.4
$ws.Range("K77").Value = 24912
$ws.Range("M77").Value = -20544

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 3001.8
$ws.Range("I25").Value = 1954.5
$ws.Range("K25").Value = 1954.5
$ws.Range("M25").Value = -1719.5
$ws.Range("H82").Value = 10714.3
$ws.Range("I82").Value = 10714.3
$ws.Range("K82").Value = 10714.3
$ws.Range("M82").Value = -10331.3
$ws.Range("H85").Value = 10714.3
$ws.Range("I85").Value = 10714.3
$ws.Range("K85").Value = 10714.3
$ws.Range("M85").Value = -9388.299999999999
$ws.Range("H86").Value = 3149.6191
$ws.Range("I86").Value = 1352.4546
$ws.Range("J86").Value = 5126.5
$ws.Range("K86").Value = 1352.4546
$ws.Range("L86").Value = 5126.5
$ws.Range("M86").Value = -229.4546
$ws.Range("N86").Value = -7372.5
$ws.Range("H89").Value = 3149.6191
$ws.Range("I89").Value = 1352.4546
$ws.Range("J89").Value = 5126.5
$ws.Range("K89").Value = 6762.273
$ws.Range("L89").Value = 25632.5
$ws.Range("M89").Value = -1146.273
$ws.Range("N89").Value = -36864.5
$ws.Range("H99").Value = 4510.5
$ws.Range("I99").Value = 4510.5
$ws.Range("K99").Value = 4510.5
$ws.Range("M99").Value = -3012.5
$ws.Range("H100").Value = 14360.75
$ws.Range("J100").Value = 14360.75
$ws.Range("L100").Value = 14360.75
$ws.Range("N100").Value = -16524.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1520.45
$ws.Range("I22").Value = 912.9091
$ws.Range("J22").Value = 2263
$ws.Range("K22").Value = 912.9091
$ws.Range("L22").Value = 2263
$ws.Range("M22").Value = -562.9091
$ws.Range("N22").Value = -2963
$ws.Range("H31").Value = 6134.6665
$ws.Range("J31").Value = 8425
$ws.Range("L31").Value = 8425
$ws.Range("N31").Value = -9015
$ws.Range("H34").Value = 6134.6665
$ws.Range("J34").Value = 8425
$ws.Range("L34").Value = 8425
$ws.Range("N34").Value = -8829
$ws.Range("H51").Value = 57972.5
$ws.Range("J51").Value = 57972.5
$ws.Range("L51").Value = 57972.5
$ws.Range("N51").Value = -59444.5
$ws.Range("H61").Value = 57972.5
$ws.Range("J61").Value = 57972.5
$ws.Range("L61").Value = 57972.5
$ws.Range("N61").Value = -58668.5
$ws.Range("H106").Value = 28064.2
$ws.Range("J106").Value = 28064.2
$ws.Range("L106").Value = 28064.2
$ws.Range("N106").Value = -30588.2
$ws.Range("H134").Value = 1876.6428
$ws.Range("I134").Value = 1905.7693
$ws.Range("J134").Value = 1498
$ws.Range("K134").Value = 5717.3079
$ws.Range("L134").Value = 4494
$ws.Range("M134").Value = -3182.3079
$ws.Range("N134").Value = -9564

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 687497.25
$ws.Range("I128").Value = 687497.25
$ws.Range("K128").Value = 2062491.75
$ws.Range("M128").Value = -2057511.75
$ws.Range("H138").Value = 2888.4443
$ws.Range("I138").Value = 2416
$ws.Range("J138").Value = 3833.3333
$ws.Range("K138").Value = 7248
$ws.Range("L138").Value = 11499.9999
$ws.Range("M138").Value = -2108
$ws.Range("N138").Value = -21779.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 5214.2856
$ws.Range("J92").Value = 6041.6665
$ws.Range("L92").Value = 6041.6665
$ws.Range("N92").Value = -9785.666499999999
$ws.Range("H122").Value = 836274.5
$ws.Range("I122").Value = 1252662.2
$ws.Range("J122").Value = 3499
$ws.Range("K122").Value = 3757986.6
$ws.Range("L122").Value = 10497
$ws.Range("M122").Value = -3755536.6
$ws.Range("N122").Value = -15397

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 459.53845
$ws.Range("I22").Value = 410
$ws.Range("K22").Value = 410
$ws.Range("M22").Value = -115
$ws.Range("H27").Value = 459.53845
$ws.Range("I27").Value = 410
$ws.Range("K27").Value = 410
$ws.Range("M27").Value = -303
$ws.Range("H46").Value = 4959.654
$ws.Range("I46").Value = 3529.6365
$ws.Range("J46").Value = 6008.3335
$ws.Range("K46").Value = 3529.6365
$ws.Range("L46").Value = 6008.3335
$ws.Range("M46").Value = -3341.6365
$ws.Range("N46").Value = -6384.3335
$ws.Range("H141").Value = 206019.75
$ws.Range("J141").Value = 206019.75
$ws.Range("L141").Value = 206019.75
$ws.Range("N141").Value = -216379.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 13333.333
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 13333.333
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 13333.333
$ws.Range("N15").Value = -13909.333
$ws.Range("M15").ClearContents()
$ws.Range("H81").Value = 1293.75
$ws.Range("I81").Value = 764.2857
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 1528.5714
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -467.5714
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 1293.75
$ws.Range("I84").Value = 764.2857
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 7642.857
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -2338.857
$ws.Range("N84").Value = -60608
$ws.Range("H97").Value = 23000
$ws.Range("J97").Value = 23000
$ws.Range("L97").Value = 23000
$ws.Range("N97").Value = -24982
